# Update vm_pu.xlsx results: slack bus voltage setpoint changed from 1.05 to 1.02 p.u.
# (case with 380 kV done) - rewrites the computed per-unit voltage magnitudes
# for rows 2-25 (buses 0-23), columns B-F and I-N.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.021616792383537
$ws.Cells.Item(2, 4).Value = 1.025696095741601
$ws.Cells.Item(2, 5).Value = 1.031928447658045
$ws.Cells.Item(2, 6).Value = 1.041801997353794
$ws.Cells.Item(2, 9).Value = 1.027832697688146
$ws.Cells.Item(2, 10).Value = 1.02680725291164
$ws.Cells.Item(2, 11).Value = 1.028520911784655
$ws.Cells.Item(2, 12).Value = 1.034735149187548
$ws.Cells.Item(2, 13).Value = 1.04458048283546
$ws.Cells.Item(2, 14).Value = 1.012994564615598

# Row 3
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.022478545607231
$ws.Cells.Item(3, 4).Value = 1.026291006740529
$ws.Cells.Item(3, 5).Value = 1.032768085685541
$ws.Cells.Item(3, 6).Value = 1.042894080136704
$ws.Cells.Item(3, 9).Value = 1.027915529060431
$ws.Cells.Item(3, 10).Value = 1.027307158528469
$ws.Cells.Item(3, 11).Value = 1.028923756964794
$ws.Cells.Item(3, 12).Value = 1.035383370296776
$ws.Cells.Item(3, 13).Value = 1.04548252944582
$ws.Cells.Item(3, 14).Value = 1.013161640491678

# Row 4
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.023036663794276
$ws.Cells.Item(4, 4).Value = 1.026676104288073
$ws.Cells.Item(4, 5).Value = 1.033312285640823
$ws.Cells.Item(4, 6).Value = 1.043602146304488
$ws.Cells.Item(4, 9).Value = 1.027967788805203
$ws.Cells.Item(4, 10).Value = 1.027630510290402
$ws.Cells.Item(4, 11).Value = 1.029183879117328
$ws.Cells.Item(4, 12).Value = 1.035803066348813
$ws.Cells.Item(4, 13).Value = 1.046067036761177
$ws.Cells.Item(4, 14).Value = 1.013269661216159

# Row 5
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.023271416552431
$ws.Cells.Item(5, 4).Value = 1.0268380336289
$ws.Cells.Item(5, 5).Value = 1.033541280854657
$ws.Cells.Item(5, 6).Value = 1.043900154979096
$ws.Cells.Item(5, 9).Value = 1.027989437951995
$ws.Cells.Item(5, 10).Value = 1.027766417418158
$ws.Cells.Item(5, 11).Value = 1.029293102681622
$ws.Cells.Item(5, 12).Value = 1.035979566145552
$ws.Cells.Item(5, 13).Value = 1.046312959433943
$ws.Cells.Item(5, 14).Value = 1.013315051523349

# Row 6
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.02331083956457
$ws.Cells.Item(6, 4).Value = 1.026865224242584
$ws.Cells.Item(6, 5).Value = 1.033579742639924
$ws.Cells.Item(6, 6).Value = 1.043950211711486
$ws.Cells.Item(6, 9).Value = 1.027993054106661
$ws.Cells.Item(6, 10).Value = 1.027789235035051
$ws.Cells.Item(6, 11).Value = 1.029311434036576
$ws.Cells.Item(6, 12).Value = 1.036009204698329
$ws.Cells.Item(6, 13).Value = 1.046354262379904
$ws.Cells.Item(6, 14).Value = 1.013322671476812

# Row 7
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.023039800102752
$ws.Cells.Item(7, 4).Value = 1.026678267863047
$ws.Cells.Item(7, 5).Value = 1.033315344648757
$ws.Cells.Item(7, 6).Value = 1.043606126985759
$ws.Cells.Item(7, 9).Value = 1.027968079343637
$ws.Cells.Item(7, 10).Value = 1.027632326407147
$ws.Cells.Item(7, 11).Value = 1.029185339087248
$ws.Cells.Item(7, 12).Value = 1.035805424515078
$ws.Cells.Item(7, 13).Value = 1.046070322024105
$ws.Cells.Item(7, 14).Value = 1.013270267809073

# Row 8
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.021907920308049
$ws.Cells.Item(8, 4).Value = 1.025897116569366
$ws.Cells.Item(8, 5).Value = 1.032212020460882
$ws.Cells.Item(8, 6).Value = 1.042170777993028
$ws.Cells.Item(8, 9).Value = 1.027860967553102
$ws.Cells.Item(8, 10).Value = 1.026976222676243
$ws.Cells.Item(8, 11).Value = 1.028657167471569
$ws.Cells.Item(8, 12).Value = 1.034954165336
$ws.Cells.Item(8, 13).Value = 1.044885162617846
$ws.Cells.Item(8, 14).Value = 1.01305104677706

# Row 9
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.01991733895273
$ws.Cells.Item(9, 4).Value = 1.024521853055343
$ws.Cells.Item(9, 5).Value = 1.030274763236917
$ws.Cells.Item(9, 6).Value = 1.039652406045211
$ws.Cells.Item(9, 9).Value = 1.027662004447214
$ws.Cells.Item(9, 10).Value = 1.02581920818236
$ws.Cells.Item(9, 11).Value = 1.02772233434898
$ws.Cells.Item(9, 12).Value = 1.033456131285591
$ws.Cells.Item(9, 13).Value = 1.042803115058044
$ws.Cells.Item(9, 14).Value = 1.012664091254353

# Row 10
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.018593008703463
$ws.Cells.Item(10, 4).Value = 1.023605933053521
$ws.Cells.Item(10, 5).Value = 1.028988010033048
$ws.Cells.Item(10, 6).Value = 1.03798089033994
$ws.Cells.Item(10, 9).Value = 1.027522527952191
$ws.Cells.Item(10, 10).Value = 1.025047340361202
$ws.Cells.Item(10, 11).Value = 1.027096406596848
$ws.Cells.Item(10, 12).Value = 1.032458850127937
$ws.Cells.Item(10, 13).Value = 1.041419426324265
$ws.Cells.Item(10, 14).Value = 1.012405701359792

# Row 11
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.018020221405071
$ws.Cells.Item(11, 4).Value = 1.023209569159728
$ws.Cells.Item(11, 5).Value = 1.028431976471061
$ws.Cells.Item(11, 6).Value = 1.03725887556717
$ws.Cells.Item(11, 9).Value = 1.027460519578512
$ws.Cells.Item(11, 10).Value = 1.024713001998317
$ws.Cells.Item(11, 11).Value = 1.026824745758826
$ws.Cells.Item(11, 12).Value = 1.032027364574016
$ws.Cells.Item(11, 13).Value = 1.040821318526489
$ws.Cells.Item(11, 14).Value = 1.012293721250572

# Row 12
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.017807562600518
$ws.Cells.Item(12, 4).Value = 1.023062379019993
$ws.Cells.Item(12, 5).Value = 1.028225613399036
$ws.Cells.Item(12, 6).Value = 1.036990953118949
$ws.Cells.Item(12, 9).Value = 1.027437244987849
$ws.Cells.Item(12, 10).Value = 1.024588797704552
$ws.Cells.Item(12, 11).Value = 1.026723745334904
$ws.Cells.Item(12, 12).Value = 1.031867144398527
$ws.Cells.Item(12, 13).Value = 1.040599311838992
$ws.Cells.Item(12, 14).Value = 1.012252112895999

# Row 13
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.01785317409466
$ws.Cells.Item(13, 4).Value = 1.023093950100774
$ws.Cells.Item(13, 5).Value = 1.028269871147447
$ws.Cells.Item(13, 6).Value = 1.037048411338108
$ws.Cells.Item(13, 9).Value = 1.027442248404534
$ws.Cells.Item(13, 10).Value = 1.024615440660449
$ws.Cells.Item(13, 11).Value = 1.026745414480895
$ws.Cells.Item(13, 12).Value = 1.031901509763981
$ws.Cells.Item(13, 13).Value = 1.040646925895897
$ws.Cells.Item(13, 14).Value = 1.012261038654362

# Row 14
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.018002640916579
$ws.Cells.Item(14, 4).Value = 1.023197401611203
$ws.Cells.Item(14, 5).Value = 1.02841491490697
$ws.Cells.Item(14, 6).Value = 1.037236723581129
$ws.Cells.Item(14, 9).Value = 1.027458600625291
$ws.Cells.Item(14, 10).Value = 1.024702735555736
$ws.Cells.Item(14, 11).Value = 1.026816398932825
$ws.Cells.Item(14, 12).Value = 1.032014119641308
$ws.Cells.Item(14, 13).Value = 1.040802964171377
$ws.Cells.Item(14, 14).Value = 1.012290282173369

# Row 15
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.018094745634821
$ws.Cells.Item(15, 4).Value = 1.023261146461776
$ws.Cells.Item(15, 5).Value = 1.028504304061304
$ws.Cells.Item(15, 6).Value = 1.037352784250339
$ws.Cells.Item(15, 9).Value = 1.027468643728984
$ws.Cells.Item(15, 10).Value = 1.024756518715009
$ws.Cells.Item(15, 11).Value = 1.026860122442272
$ws.Cells.Item(15, 12).Value = 1.032083509329661
$ws.Cells.Item(15, 13).Value = 1.040899125365033
$ws.Cells.Item(15, 14).Value = 1.01230829823334

# Row 16
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.018631036649519
$ws.Cells.Item(16, 4).Value = 1.023632243524825
$ws.Cells.Item(16, 5).Value = 1.029024936269165
$ws.Cells.Item(16, 6).Value = 1.038028845318304
$ws.Cells.Item(16, 9).Value = 1.027526609285791
$ws.Cells.Item(16, 10).Value = 1.025069526975546
$ws.Cells.Item(16, 11).Value = 1.027114422657494
$ws.Cells.Item(16, 12).Value = 1.032487493741775
$ws.Cells.Item(16, 13).Value = 1.041459142748039
$ws.Cells.Item(16, 14).Value = 1.012413131131553

# Row 17
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.01896761448342
$ws.Cells.Item(17, 4).Value = 1.023865087145334
$ws.Cells.Item(17, 5).Value = 1.029351821059544
$ws.Cells.Item(17, 6).Value = 1.038453393284052
$ws.Cells.Item(17, 9).Value = 1.027562537812468
$ws.Cells.Item(17, 10).Value = 1.025265838882734
$ws.Cells.Item(17, 11).Value = 1.027273770613996
$ws.Cells.Item(17, 12).Value = 1.032740995319694
$ws.Cells.Item(17, 13).Value = 1.041810705716674
$ws.Cells.Item(17, 14).Value = 1.012478864734606

# Row 18
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.019163997890344
$ws.Cells.Item(18, 4).Value = 1.024000923501011
$ws.Cells.Item(18, 5).Value = 1.029542597235513
$ws.Cells.Item(18, 6).Value = 1.03870119491138
$ws.Cells.Item(18, 9).Value = 1.027583338534726
$ws.Cells.Item(18, 10).Value = 1.025380333177665
$ws.Cells.Item(18, 11).Value = 1.027366654718616
$ws.Cells.Item(18, 12).Value = 1.032888891594376
$ws.Cells.Item(18, 13).Value = 1.042015866576519
$ws.Cells.Item(18, 14).Value = 1.012517196784025

# Row 19
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.019230970246692
$ws.Cells.Item(19, 4).Value = 1.024047243972197
$ws.Cells.Item(19, 5).Value = 1.029607665563668
$ws.Cells.Item(19, 6).Value = 1.038785717658518
$ws.Cells.Item(19, 9).Value = 1.027590404602775
$ws.Cells.Item(19, 10).Value = 1.025419370855088
$ws.Cells.Item(19, 11).Value = 1.027398315439439
$ws.Cells.Item(19, 12).Value = 1.03293932598108
$ws.Cells.Item(19, 13).Value = 1.042085838076519
$ws.Cells.Item(19, 14).Value = 1.012530265441191

# Row 20
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.018931496315599
$ws.Cells.Item(20, 4).Value = 1.023840102877571
$ws.Cells.Item(20, 5).Value = 1.029316738046912
$ws.Cells.Item(20, 6).Value = 1.038407825714646
$ws.Cells.Item(20, 9).Value = 1.027558699130357
$ws.Cells.Item(20, 10).Value = 1.025244777615392
$ws.Cells.Item(20, 11).Value = 1.027256680373251
$ws.Cells.Item(20, 12).Value = 1.032713793568387
$ws.Cells.Item(20, 13).Value = 1.041772975972254
$ws.Cells.Item(20, 14).Value = 1.012471813094856

# Row 21
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.017958623894828
$ws.Cells.Item(21, 4).Value = 1.023166936683407
$ws.Cells.Item(21, 5).Value = 1.028372198360279
$ws.Cells.Item(21, 6).Value = 1.037181262967079
$ws.Cells.Item(21, 9).Value = 1.027453791977315
$ws.Cells.Item(21, 10).Value = 1.024677029821826
$ws.Cells.Item(21, 11).Value = 1.026795498352557
$ws.Cells.Item(21, 12).Value = 1.031980957373386
$ws.Cells.Item(21, 13).Value = 1.040757010440235
$ws.Cells.Item(21, 14).Value = 1.012281671067913

# Row 22
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.017347519025398
$ws.Cells.Item(22, 4).Value = 1.022743905615136
$ws.Cells.Item(22, 5).Value = 1.027779327796715
$ws.Cells.Item(22, 6).Value = 1.036411613788607
$ws.Cells.Item(22, 9).Value = 1.027386433303233
$ws.Cells.Item(22, 10).Value = 1.024319971399437
$ws.Cells.Item(22, 11).Value = 1.026504994561749
$ws.Cells.Item(22, 12).Value = 1.031520500046637
$ws.Cells.Item(22, 13).Value = 1.040119142549863
$ws.Cells.Item(22, 14).Value = 1.012162040682839

# Row 23
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.017671422003243
$ws.Cells.Item(23, 4).Value = 1.022968141428046
$ws.Cells.Item(23, 5).Value = 1.028093524529765
$ws.Cells.Item(23, 6).Value = 1.036819473057037
$ws.Cells.Item(23, 9).Value = 1.02742227386072
$ws.Cells.Item(23, 10).Value = 1.024509263270906
$ws.Cells.Item(23, 11).Value = 1.026659046951015
$ws.Cells.Item(23, 12).Value = 1.031764567818484
$ws.Cells.Item(23, 13).Value = 1.040457201748744
$ws.Cells.Item(23, 14).Value = 1.012225466513557

# Row 24
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.018947816365638
$ws.Cells.Item(24, 4).Value = 1.023851392123074
$ws.Cells.Item(24, 5).Value = 1.029332590211053
$ws.Cells.Item(24, 6).Value = 1.038428415212706
$ws.Cells.Item(24, 9).Value = 1.027560434146966
$ws.Cells.Item(24, 10).Value = 1.025254294330002
$ws.Cells.Item(24, 11).Value = 1.027264402905586
$ws.Cells.Item(24, 12).Value = 1.032726084766778
$ws.Cells.Item(24, 13).Value = 1.041790024110418
$ws.Cells.Item(24, 14).Value = 1.012474999455949

# Row 25
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.020431477549766
$ws.Cells.Item(25, 4).Value = 1.024877236290779
$ws.Cells.Item(25, 5).Value = 1.030774759382809
$ws.Cells.Item(25, 6).Value = 1.04030216680145
$ws.Cells.Item(25, 9).Value = 1.027714648421421
$ws.Cells.Item(25, 10).Value = 1.026118421359177
$ws.Cells.Item(25, 11).Value = 1.027964493078835
$ws.Cells.Item(25, 12).Value = 1.033843165771968
$ws.Cells.Item(25, 13).Value = 1.043340613951886
$ws.Cells.Item(25, 14).Value = 1.012764204253668

